# Insert a new data row at row 17 (pushes existing rows 17-130 down to 18-131)
# and populate the new row with the Macroferia Regional de Talca - Alcachofa
# "Madrigal / Primera" reading for the new date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17:17").Insert()

$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 45149
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 100112013
$ws.Range("G17").Value = "Alcachofa"
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 12000
$ws.Range("N17").Value = "`$/caja 40 unidades"
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 300
$ws.Range("Q17").Value = 40
$ws.Range("R17").Value = "Hortaliza"
